# Developer guide image update: rename "Person*" UI component shapes to
# "Member*" on the single diagram slide, and shrink the PersonListPanel
# label to fit its box.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        if ($shp.TextFrame.HasText) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "PersonListPanel") {
                $tr.Text = "MemberListPanel"
                $tr.Font.Size = 9
            }
            elseif ($tr.Text -eq "PersonCard") {
                $tr.Text = "MemberCard"
            }
        }
    }
}
